$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data values (B:E), with G recomputed as the row sum (B+C+D+E+F)
$data = @{
    2 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732)
    3 = @(0.127881588408715,  0.3127903958511391, 3.900430680208489, 0.496779210170732)
    4 = @(0.3048080303191223, 1.667794583268128, 0.8054896365839992, 8.660232485948974)
    5 = @(3.230985683306322,  1.667794583268128, 0.1575252929769615, 0.496779210170732)
    6 = @(1.459612070389937,  1.667794583268128, 0.8054896365839992, 0.496779210170732)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value2 = $vals[0]
    $ws.Cells.Item($row, 3).Value2 = $vals[1]
    $ws.Cells.Item($row, 4).Value2 = $vals[2]
    $ws.Cells.Item($row, 5).Value2 = $vals[3]

    $fVal = $ws.Cells.Item($row, 6).Value2
    $sum = $vals[0] + $vals[1] + $vals[2] + $vals[3] + $fVal
    $ws.Cells.Item($row, 7).Value2 = $sum
}
